$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(44946, "Candy White", "Primera", 20, 320000, 330000, 325000, "`$/bins (420 kilos)", "Región de O'Higgins", 774, 420)
    ,@(44946, "Candy White", "Segunda", 10, 290000, 300000, 295000, "`$/bins (420 kilos)", "Región de O'Higgins", 702, 420)
    ,@(44946, "Magique", "Especial", 10, 510000, 520000, 515000, "`$/bins (420 kilos)", "Región de O'Higgins", 1226, 420)
    ,@(44946, "Magique", "Primera", 10, 430000, 440000, 435000, "`$/bins (420 kilos)", "Región de O'Higgins", 1036, 420)
    ,@(44224, "Nectar Crest", "Especial", 140, 17500, 18000, 17750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1109, 16)
    ,@(44224, "Nectar Crest", "Primera", 240, 15500, 16000, 15750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 984, 16)
    ,@(44224, "Nectar Crest", "Segunda", 200, 13500, 14000, 13750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 859, 16)
    ,@(44224, "Venus", "Especial", 140, 18500, 19000, 18750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1172, 16)
    ,@(44224, "Venus", "Primera", 240, 16500, 17000, 16750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1047, 16)
    ,@(44224, "Venus", "Segunda", 200, 14500, 15000, 14750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 922, 16)
    ,@(44637, "August Red", "Especial", 16, 450000, 460000, 455000, "`$/bins (420 kilos)", "Región de O'Higgins", 1083, 420)
    ,@(44637, "August Red", "Primera", 20, 420000, 430000, 425000, "`$/bins (420 kilos)", "Región de O'Higgins", 1012, 420)
    ,@(44637, "August Red", "Segunda", 20, 380000, 390000, 385000, "`$/bins (420 kilos)", "Región de O'Higgins", 917, 420)
    ,@(44239, "Venus", "Especial", 100, 19000, 19500, 19250, "`$/caja 16 kilos empedrada", "Región Metropolitana", 1203, 16)
    ,@(44239, "Venus", "Primera", 160, 17000, 17500, 17250, "`$/caja 16 kilos empedrada", "Región Metropolitana", 1078, 16)
    ,@(44239, "Venus", "Segunda", 160, 15000, 15500, 15250, "`$/caja 16 kilos empedrada", "Región Metropolitana", 953, 16)
    ,@(44616, "Artic Snow", "Especial", 16, 335000, 340000, 337500, "`$/bins (420 kilos)", "Región de O'Higgins", 804, 420)
    ,@(44616, "Artic Snow", "Primera", 20, 305000, 310000, 307500, "`$/bins (420 kilos)", "Región de O'Higgins", 732, 420)
    ,@(44616, "August Red", "Primera", 10, 315000, 320000, 317500, "`$/bins (420 kilos)", "Región de O'Higgins", 756, 420)
    ,@(44616, "August Red", "Segunda", 10, 275000, 280000, 277500, "`$/bins (420 kilos)", "Región de O'Higgins", 661, 420)
    ,@(44616, "Venus", "Especial", 28, 325000, 330000, 326786, "`$/bins (420 kilos)", "Región de O'Higgins", 778, 420)
    ,@(44616, "Venus", "Primera", 20, 305000, 310000, 307500, "`$/bins (420 kilos)", "Región de O'Higgins", 732, 420)
    ,@(44616, "Venus", "Segunda", 20, 255000, 260000, 257500, "`$/bins (420 kilos)", "Región de O'Higgins", 613, 420)
    ,@(44650, "Artic Snow", "Especial", 16, 410000, 420000, 415000, "`$/bins (420 kilos)", "Región de O'Higgins", 988, 420)
    ,@(44650, "Artic Snow", "Primera", 14, 370000, 380000, 375000, "`$/bins (420 kilos)", "Región de O'Higgins", 893, 420)
    ,@(44650, "Artic Snow", "Segunda", 10, 350000, 360000, 355000, "`$/bins (420 kilos)", "Región de O'Higgins", 845, 420)
    ,@(44581, "Ruby Diamond", "Especial", 20, 385000, 390000, 387500, "`$/bins (420 kilos)", "Región de O'Higgins", 923, 420)
    ,@(44581, "Ruby Diamond", "Primera", 20, 335000, 340000, 337500, "`$/bins (420 kilos)", "Región de O'Higgins", 804, 420)
    ,@(44581, "Super Queen", "Especial", 20, 355000, 360000, 357500, "`$/bins (420 kilos)", "Región de O'Higgins", 851, 420)
    ,@(44581, "Super Queen", "Primera", 20, 325000, 330000, 327500, "`$/bins (420 kilos)", "Región de O'Higgins", 780, 420)
    ,@(44581, "Super Queen", "Segunda", 20, 305000, 310000, 307500, "`$/bins (420 kilos)", "Región de O'Higgins", 732, 420)
    ,@(44631, "Artic Snow", "Primera", 16, 330000, 335000, 332500, "`$/bins (420 kilos)", "Región de O'Higgins", 792, 420)
    ,@(44631, "Artic Snow", "Segunda", 16, 290000, 295000, 292500, "`$/bins (420 kilos)", "Región de O'Higgins", 696, 420)
    ,@(44643, "Artic Snow", "Primera", 16, 430000, 440000, 435000, "`$/bins (420 kilos)", "Región de O'Higgins", 1036, 420)
    ,@(44643, "Artic Snow", "Segunda", 12, 400000, 410000, 405000, "`$/bins (420 kilos)", "Región de O'Higgins", 964, 420)
    ,@(44643, "August Red", "Primera", 18, 430000, 440000, 435000, "`$/bins (420 kilos)", "Región de O'Higgins", 1036, 420)
    ,@(44643, "August Red", "Segunda", 12, 370000, 380000, 375000, "`$/bins (420 kilos)", "Región de O'Higgins", 893, 420)
    ,@(44271, "June Pearl", "Especial", 240, 22000, 23000, 22500, "`$/caja 18 kilos empedrada", "Región de O'Higgins", 1250, 18)
    ,@(44271, "June Pearl", "Primera", 160, 20000, 21000, 20500, "`$/caja 18 kilos empedrada", "Región de O'Higgins", 1139, 18)
    ,@(44908, "Artic Star", "Primera", 20, 350000, 360000, 355000, "`$/bins (420 kilos)", "Región de O'Higgins", 845, 420)
    ,@(44908, "Artic Star", "Segunda", 16, 320000, 330000, 325000, "`$/bins (420 kilos)", "Región de O'Higgins", 774, 420)
    ,@(44908, "Super Queen", "Especial", 20, 450000, 460000, 455000, "`$/bins (420 kilos)", "Región de O'Higgins", 1083, 420)
    ,@(44908, "Super Queen", "Primera", 16, 400000, 410000, 405000, "`$/bins (420 kilos)", "Región de O'Higgins", 964, 420)
    ,@(44259, "August Red", "Especial", 200, 19500, 20000, 19750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1234, 16)
    ,@(44259, "August Red", "Primera", 200, 17500, 18000, 17750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1109, 16)
    ,@(44259, "August Red", "Segunda", 200, 15500, 16000, 15750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 984, 16)
    ,@(44252, "August Red", "Especial", 160, 17500, 18000, 17750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1109, 16)
    ,@(44252, "August Red", "Primera", 200, 15500, 16000, 15750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 984, 16)
    ,@(44252, "August Red", "Segunda", 240, 13500, 14000, 13750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 859, 16)
    ,@(44252, "June Pearl", "Especial", 160, 17500, 18000, 17750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1109, 16)
    ,@(44252, "June Pearl", "Primera", 240, 15500, 16000, 15750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 984, 16)
    ,@(44252, "June Pearl", "Segunda", 240, 13500, 14000, 13750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 859, 16)
    ,@(44235, "June Pearl", "Especial", 100, 19500, 20000, 19750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1234, 16)
    ,@(44235, "June Pearl", "Primera", 240, 17500, 18000, 17750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1109, 16)
    ,@(44235, "June Pearl", "Segunda", 160, 14500, 15000, 14750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 922, 16)
    ,@(44235, "Venus", "Especial", 140, 19500, 20000, 19750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1234, 16)
    ,@(44235, "Venus", "Primera", 160, 17500, 18000, 17750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1109, 16)
    ,@(44235, "Venus", "Segunda", 100, 15500, 16000, 15750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 984, 16)
    ,@(44539, "Artic Star", "Especial", 10, 485000, 490000, 487500, "`$/bins (420 kilos)", "Región de O'Higgins", 1161, 420)
    ,@(44539, "Artic Star", "Primera", 20, 455000, 460000, 457500, "`$/bins (420 kilos)", "Región de O'Higgins", 1089, 420)
    ,@(44539, "Artic Star", "Segunda", 16, 425000, 430000, 427500, "`$/bins (420 kilos)", "Región de O'Higgins", 1018, 420)
    ,@(44281, "August Red", "Primera", 20, 295000, 300000, 297500, "`$/bins (420 kilos)", "Región de O'Higgins", 708, 420)
    ,@(44281, "August Red", "Segunda", 20, 265000, 270000, 267500, "`$/bins (420 kilos)", "Región de O'Higgins", 637, 420)
    ,@(44214, "Red Diamond", "Especial", 200, 20500, 21000, 20750, "`$/caja 17 kilos empedrada", "Región de O'Higgins", 1221, 17)
    ,@(44214, "Red Diamond", "Primera", 200, 18500, 19000, 18750, "`$/caja 17 kilos empedrada", "Región de O'Higgins", 1103, 17)
    ,@(44214, "Red Diamond", "Segunda", 200, 16500, 17000, 16750, "`$/caja 17 kilos empedrada", "Región de O'Higgins", 985, 17)
    ,@(44214, "Venus", "Especial", 240, 20500, 21000, 20750, "`$/caja 17 kilos empedrada", "Región de O'Higgins", 1221, 17)
    ,@(44214, "Venus", "Primera", 200, 18500, 19000, 18750, "`$/caja 17 kilos empedrada", "Región de O'Higgins", 1103, 17)
    ,@(44214, "Venus", "Segunda", 160, 16500, 17000, 16750, "`$/caja 17 kilos empedrada", "Región de O'Higgins", 985, 17)
    ,@(44592, "June Pearl", "Especial", 16, 405000, 410000, 407500, "`$/bins (420 kilos)", "Región de O'Higgins", 970, 420)
    ,@(44592, "June Pearl", "Primera", 16, 355000, 360000, 357500, "`$/bins (420 kilos)", "Región de O'Higgins", 851, 420)
    ,@(44592, "Venus", "Especial", 16, 375000, 380000, 377500, "`$/bins (420 kilos)", "Región de O'Higgins", 899, 420)
    ,@(44592, "Venus", "Primera", 16, 335000, 340000, 337500, "`$/bins (420 kilos)", "Región de O'Higgins", 804, 420)
    ,@(44225, "Nectar Crest", "Especial", 100, 17500, 18000, 17750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1109, 16)
    ,@(44225, "Nectar Crest", "Primera", 200, 15500, 16000, 15750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 984, 16)
    ,@(44225, "Nectar Crest", "Segunda", 200, 13500, 14000, 13750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 859, 16)
    ,@(44225, "Venus", "Especial", 100, 18500, 19000, 18750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1172, 16)
    ,@(44225, "Venus", "Primera", 200, 16500, 17000, 16750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 1047, 16)
    ,@(44225, "Venus", "Segunda", 200, 14500, 15000, 14750, "`$/caja 16 kilos empedrada", "Región de O'Higgins", 922, 16)
    ,@(44636, "Artic Snow", "Especial", 20, 390000, 400000, 395000, "`$/bins (420 kilos)", "Región de O'Higgins", 940, 420)
    ,@(44595, "June Pearl", "Especial", 20, 400000, 405000, 402500, "`$/bins (420 kilos)", "Región de O'Higgins", 958, 420)
    ,@(44595, "June Pearl", "Primera", 20, 350000, 360000, 355000, "`$/bins (420 kilos)", "Región de O'Higgins", 845, 420)
    ,@(44595, "June Pearl", "Segunda", 20, 325000, 330000, 327500, "`$/bins (420 kilos)", "Región de O'Higgins", 780, 420)
    ,@(44595, "Venus", "Especial", 20, 370000, 375000, 372500, "`$/bins (420 kilos)", "Región de O'Higgins", 887, 420)
    ,@(44595, "Venus", "Primera", 28, 330000, 335000, 331786, "`$/bins (420 kilos)", "Región de O'Higgins", 790, 420)
)

# Constant columns shared by every data row in this sheet
$constMercadoId = 8
$constMercado = "Terminal La Palmera de La Serena"
$constRegion = "Coquimbo"
$constCodreg = 4
$constTipo = "Fruta"
$constProductoId = 100103
$constProducto = "Frutos de hueso (carozo)"
$constCategoriaId = 100103006
$constCategoria = "Nectarín"

$startRow = 354
$lastRowBefore = 434
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    if ($row -gt $lastRowBefore) {
        # Brand-new row: populate the constant identifying columns too
        $ws.Cells.Item($row, 1).Value = $constMercadoId
        $ws.Cells.Item($row, 2).Value = $constMercado
        $ws.Cells.Item($row, 3).Value = $constRegion
        $ws.Cells.Item($row, 5).Value = $constCodreg
        $ws.Cells.Item($row, 6).Value = $constTipo
        $ws.Cells.Item($row, 7).Value = $constProductoId
        $ws.Cells.Item($row, 8).Value = $constProducto
        $ws.Cells.Item($row, 9).Value = $constCategoriaId
        $ws.Cells.Item($row, 10).Value = $constCategoria
        # Match the date-column number format used throughout column D
        $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 11).Value = $vals[1]
    $ws.Cells.Item($row, 12).Value = $vals[2]
    $ws.Cells.Item($row, 13).Value = $vals[3]
    $ws.Cells.Item($row, 14).Value = $vals[4]
    $ws.Cells.Item($row, 15).Value = $vals[5]
    $ws.Cells.Item($row, 16).Value = $vals[6]
    $ws.Cells.Item($row, 17).Value = $vals[7]
    $ws.Cells.Item($row, 18).Value = $vals[8]
    $ws.Cells.Item($row, 19).Value = $vals[9]
    $ws.Cells.Item($row, 20).Value = $vals[10]
}

$ws.Range("A1").Select()